# Update computed line-power results (pl_mw) for the 380 kV case.
# Columns: B,I,L,N remain 0; C..H,J,K,M,O receive new values for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.04904822163325662
$ws.Cells.Item(2, 4).Value = 0.2150398836567859
$ws.Cells.Item(2, 5).Value = 0.1577557748486988
$ws.Cells.Item(2, 6).Value = 1.098088969301038
$ws.Cells.Item(2, 7).Value = 0.5478667889177586
$ws.Cells.Item(2, 8).Value = 0.6919128742687448
$ws.Cells.Item(2, 10).Value = 0.1587236157017031
$ws.Cells.Item(2, 11).Value = 1.264877204211814
$ws.Cells.Item(2, 13).Value = 0.3963371515211236
$ws.Cells.Item(2, 15).Value = 2.452149567775407
$ws.Cells.Item(3, 3).Value = 0.04357218035039523
$ws.Cells.Item(3, 4).Value = 0.2123531160951586
$ws.Cells.Item(3, 5).Value = 0.1580148820677998
$ws.Cells.Item(3, 6).Value = 1.107628383791486
$ws.Cells.Item(3, 7).Value = 0.554649574212128
$ws.Cells.Item(3, 8).Value = 0.6999168003676317
$ws.Cells.Item(3, 10).Value = 0.1604476692235615
$ws.Cells.Item(3, 11).Value = 1.111272602920906
$ws.Cells.Item(3, 13).Value = 0.364008585198242
$ws.Cells.Item(3, 15).Value = 2.482981324815299
$ws.Cells.Item(4, 3).Value = 0.04022327421418481
$ws.Cells.Item(4, 4).Value = 0.2107731110297237
$ws.Cells.Item(4, 5).Value = 0.1582495810300522
$ws.Cells.Item(4, 6).Value = 1.114257136937795
$ws.Cells.Item(4, 7).Value = 0.5593307832977743
$ws.Cells.Item(4, 8).Value = 0.7052314136264073
$ws.Cells.Item(4, 10).Value = 0.1616012591945655
$ws.Cells.Item(4, 11).Value = 1.016639218084151
$ws.Cells.Item(4, 13).Value = 0.3441705635275412
$ws.Cells.Item(4, 15).Value = 2.503834580801595
$ws.Cells.Item(5, 3).Value = 0.0388619440559097
$ws.Cells.Item(5, 4).Value = 0.2101468542831881
$ws.Cells.Item(5, 5).Value = 0.1583642538953178
$ws.Cells.Item(5, 6).Value = 1.117152235090536
$ws.Cells.Item(5, 7).Value = 0.561367998431507
$ws.Cells.Item(5, 8).Value = 0.7074977416082291
$ws.Cells.Item(5, 10).Value = 0.1620952355367749
$ws.Cells.Item(5, 11).Value = 0.9779976586054886
$ws.Cells.Item(5, 13).Value = 0.3360899392158743
$ws.Cells.Item(5, 15).Value = 2.512815109786331
$ws.Cells.Item(6, 3).Value = 0.03863610030231257
$ws.Cells.Item(6, 4).Value = 0.2100439308608841
$ws.Cells.Item(6, 5).Value = 0.1583844452415111
$ws.Cells.Item(6, 6).Value = 1.117644665458343
$ws.Cells.Item(6, 7).Value = 0.5617140940216601
$ws.Cells.Item(6, 8).Value = 0.7078801370699139
$ws.Cells.Item(6, 10).Value = 0.1621787017868748
$ws.Cells.Item(6, 11).Value = 0.9715766438862659
$ws.Cells.Item(6, 13).Value = 0.3347483872520982
$ws.Cells.Item(6, 15).Value = 2.514335448994615
$ws.Cells.Item(7, 3).Value = 0.04020490116032249
$ws.Cells.Item(7, 4).Value = 0.2107645937025637
$ws.Cells.Item(7, 5).Value = 0.1582510504662409
$ws.Cells.Item(7, 6).Value = 1.114295396679452
$ws.Cells.Item(7, 7).Value = 0.5593577336190805
$ws.Cells.Item(7, 8).Value = 0.7052615709365142
$ws.Cells.Item(7, 10).Value = 0.1616078244641592
$ws.Cells.Item(7, 11).Value = 1.016118395500115
$ws.Cells.Item(7, 13).Value = 0.3440615703729506
$ws.Cells.Item(7, 15).Value = 2.503953742097536
$ws.Cells.Item(8, 3).Value = 0.04715730191674083
$ws.Cells.Item(8, 4).Value = 0.2140990743388187
$ws.Cells.Item(8, 5).Value = 0.1578294342764046
$ws.Cells.Item(8, 6).Value = 1.101217925346134
$ws.Cells.Item(8, 7).Value = 0.5500980989123718
$ws.Cells.Item(8, 8).Value = 0.6945895453361786
$ws.Cells.Item(8, 10).Value = 0.159298341406604
$ws.Cells.Item(8, 11).Value = 1.211982202247157
$ws.Cells.Item(8, 13).Value = 0.3851881218608995
$ws.Cells.Item(8, 15).Value = 2.462380905176076
$ws.Cells.Item(9, 3).Value = 0.06089776901349353
$ws.Cells.Item(9, 4).Value = 0.2211875795123746
$ws.Cells.Item(9, 5).Value = 0.1576018437688269
$ws.Cells.Item(9, 6).Value = 1.081703099948257
$ws.Cells.Item(9, 7).Value = 0.5360525144575945
$ws.Cells.Item(9, 8).Value = 0.6768389332436158
$ws.Cells.Item(9, 10).Value = 0.1555239397082175
$ws.Cells.Item(9, 11).Value = 1.59343712896964
$ws.Cells.Item(9, 13).Value = 0.4659111152231787
$ws.Cells.Item(9, 15).Value = 2.396144412355426
$ws.Cells.Item(10, 3).Value = 0.07105968062690238
$ws.Cells.Item(10, 4).Value = 0.2267268712667772
$ws.Cells.Item(10, 5).Value = 0.1577991947337658
$ws.Cells.Item(10, 6).Value = 1.071115004803048
$ws.Cells.Item(10, 7).Value = 0.5282597451557507
$ws.Cells.Item(10, 8).Value = 0.6657371925718962
$ws.Cells.Item(10, 10).Value = 0.1532116730378412
$ws.Cells.Item(10, 11).Value = 1.871982938775091
$ws.Cells.Item(10, 13).Value = 0.5252425246772816
$ws.Cells.Item(10, 15).Value = 2.356847562035739
$ws.Cells.Item(11, 3).Value = 0.07569751243225653
$ws.Cells.Item(11, 4).Value = 0.2293180620887938
$ws.Cells.Item(11, 5).Value = 0.1579679976980302
$ws.Cells.Item(11, 6).Value = 1.067115070208644
$ws.Cells.Item(11, 7).Value = 0.5252673879053447
$ws.Cells.Item(11, 8).Value = 0.6611084905914879
$ws.Cells.Item(11, 10).Value = 0.1522599998518999
$ws.Cells.Item(11, 11).Value = 1.99830960368962
$ws.Cells.Item(11, 13).Value = 0.5522351947000033
$ws.Cells.Item(11, 15).Value = 2.34101455040863
$ws.Cells.Item(12, 3).Value = 0.07745592301635895
$ws.Cells.Item(12, 4).Value = 0.2303094613213261
$ws.Cells.Item(12, 5).Value = 0.1580432676107222
$ws.Cells.Item(12, 6).Value = 1.065718041412481
$ws.Cells.Item(12, 7).Value = 0.5242140604122056
$ws.Cells.Item(12, 8).Value = 0.6594163941489626
$ws.Cells.Item(12, 10).Value = 0.1519140488595028
$ws.Cells.Item(12, 11).Value = 2.04608858929106
$ws.Cells.Item(12, 13).Value = 0.5624565166892523
$ws.Cells.Item(12, 15).Value = 2.33531366377585
$ws.Cells.Item(13, 3).Value = 0.07707712175698589
$ws.Cells.Item(13, 4).Value = 0.230095494758416
$ws.Cells.Item(13, 5).Value = 0.158026552444241
$ws.Cells.Item(13, 6).Value = 1.066013680014429
$ws.Cells.Item(13, 7).Value = 0.5244373580929675
$ws.Cells.Item(13, 8).Value = 0.6597781168642598
$ws.Cells.Item(13, 10).Value = 0.1519879136335192
$ws.Cells.Item(13, 11).Value = 2.035801161549784
$ws.Cells.Item(13, 13).Value = 0.5602551908890945
$ws.Cells.Item(13, 15).Value = 2.336528329115595
$ws.Cells.Item(14, 3).Value = 0.07584213461770162
$ws.Cells.Item(14, 4).Value = 0.2293994217452848
$ws.Cells.Item(14, 5).Value = 0.1579739628626129
$ws.Cells.Item(14, 6).Value = 1.066997776438519
$ws.Cells.Item(14, 7).Value = 0.5251791284538996
$ws.Cells.Item(14, 8).Value = 0.6609680639940052
$ws.Cells.Item(14, 10).Value = 0.1522312490025612
$ws.Cells.Item(14, 11).Value = 2.002241591484051
$ws.Cells.Item(14, 13).Value = 0.5530761166408809
$ws.Cells.Item(14, 15).Value = 2.340539623019453
$ws.Cells.Item(15, 3).Value = 0.07508595086341074
$ws.Cells.Item(15, 4).Value = 0.2289743789804675
$ws.Cells.Item(15, 5).Value = 0.1579432275782331
$ws.Cells.Item(15, 6).Value = 1.06761589338786
$ws.Cells.Item(15, 7).Value = 0.5256438885627972
$ws.Cells.Item(15, 8).Value = 0.6617048471901938
$ws.Cells.Item(15, 10).Value = 0.1523821783550616
$ws.Cells.Item(15, 11).Value = 1.981677726717976
$ws.Cells.Item(15, 13).Value = 0.5486786811282656
$ws.Cells.Item(15, 15).Value = 2.34303506581972
$ws.Cells.Item(16, 3).Value = 0.07075689088362935
$ws.Cells.Item(16, 4).Value = 0.2265589589354846
$ws.Cells.Item(16, 5).Value = 0.1577897519132776
$ws.Cells.Item(16, 6).Value = 1.071392862355232
$ws.Cells.Item(16, 7).Value = 0.5284664525289315
$ws.Cells.Item(16, 8).Value = 0.6660481780263297
$ws.Cells.Item(16, 10).Value = 0.1532758848787736
$ws.Cells.Item(16, 11).Value = 1.863719203685037
$ws.Cells.Item(16, 13).Value = 0.523478491298377
$ws.Cells.Item(16, 15).Value = 2.357923480829527
$ws.Cells.Item(17, 3).Value = 0.0681050195957198
$ws.Cells.Item(17, 4).Value = 0.2250953866446821
$ws.Cells.Item(17, 5).Value = 0.1577158288326572
$ws.Cells.Item(17, 6).Value = 1.07391923466237
$ws.Cells.Item(17, 7).Value = 0.5303397898178233
$ws.Cells.Item(17, 8).Value = 0.6688206904866618
$ws.Cells.Item(17, 10).Value = 0.1538498198213532
$ws.Cells.Item(17, 11).Value = 1.791254814902118
$ws.Cells.Item(17, 13).Value = 0.5080192223795592
$ws.Cells.Item(17, 15).Value = 2.367581029739711
$ws.Cells.Item(18, 3).Value = 0.06658115789831243
$ws.Cells.Item(18, 4).Value = 0.2242602985358104
$ws.Cells.Item(18, 5).Value = 0.1576807495932293
$ws.Cells.Item(18, 6).Value = 1.075449190011305
$ws.Cells.Item(18, 7).Value = 0.5314692776900714
$ws.Cells.Item(18, 8).Value = 0.6704550386926655
$ws.Cells.Item(18, 10).Value = 0.1541893610465763
$ws.Cells.Item(18, 11).Value = 1.749539121086116
$ws.Cells.Item(18, 13).Value = 0.4991277306812378
$ws.Cells.Item(18, 15).Value = 2.373328069871036
$ws.Cells.Item(19, 3).Value = 0.06606544950081172
$ws.Cells.Item(19, 4).Value = 0.22397870886671
$ws.Cells.Item(19, 5).Value = 0.1576701505730291
$ws.Cells.Item(19, 6).Value = 1.075980399196602
$ws.Cells.Item(19, 7).Value = 0.5318606222236539
$ws.Cells.Item(19, 8).Value = 0.6710152123568562
$ws.Cells.Item(19, 10).Value = 0.1543059427057152
$ws.Cells.Item(19, 11).Value = 1.73540879010352
$ws.Cells.Item(19, 13).Value = 0.4961172885814591
$ws.Cells.Item(19, 15).Value = 2.3753069144514
$ws.Cells.Item(20, 3).Value = 0.06838716846735338
$ws.Cells.Item(20, 4).Value = 0.2252504915547036
$ws.Cells.Item(20, 5).Value = 0.1577229282415331
$ws.Cells.Item(20, 6).Value = 1.073642342385369
$ws.Cells.Item(20, 7).Value = 0.5301349861715678
$ws.Cells.Item(20, 8).Value = 0.6685214453067729
$ws.Cells.Item(20, 10).Value = 0.1537877474530198
$ws.Cells.Item(20, 11).Value = 1.798972528658965
$ws.Cells.Item(20, 13).Value = 0.5096648641109454
$ws.Cells.Item(20, 15).Value = 2.36653306102636
$ws.Cells.Item(21, 3).Value = 0.07620482162614906
$ws.Cells.Item(21, 4).Value = 0.2296035999019637
$ws.Cells.Item(21, 5).Value = 0.1579891018473916
$ws.Cells.Item(21, 6).Value = 1.066705528412029
$ws.Cells.Item(21, 7).Value = 0.5249590836432532
$ws.Cells.Item(21, 8).Value = 0.6606168998925313
$ws.Cells.Item(21, 10).Value = 0.1521593838026831
$ws.Cells.Item(21, 11).Value = 2.012100451789593
$ws.Cells.Item(21, 13).Value = 0.5551847952243776
$ws.Cells.Item(21, 15).Value = 2.339353402660691
$ws.Cells.Item(22, 3).Value = 0.08132671992625262
$ws.Cells.Item(22, 4).Value = 0.2325078550610584
$ws.Cells.Item(22, 5).Value = 0.1582291984123749
$ws.Cells.Item(22, 6).Value = 1.062857834109536
$ws.Cells.Item(22, 7).Value = 0.5220417051451278
$ws.Cells.Item(22, 8).Value = 0.6558046203546581
$ws.Cells.Item(22, 10).Value = 0.1511792547930746
$ws.Cells.Item(22, 11).Value = 2.151051474068595
$ws.Cells.Item(22, 13).Value = 0.5849331757904963
$ws.Cells.Item(22, 15).Value = 2.323308260375256
$ws.Cells.Item(23, 3).Value = 0.07859191630960538
$ws.Cells.Item(23, 4).Value = 0.2309524069445388
$ws.Cells.Item(23, 5).Value = 0.1580950081644303
$ws.Cells.Item(23, 6).Value = 1.064848585063501
$ws.Cells.Item(23, 7).Value = 0.5235560687179657
$ws.Cells.Item(23, 8).Value = 0.6583406268412801
$ws.Cells.Item(23, 10).Value = 0.1516946662743628
$ws.Cells.Item(23, 11).Value = 2.076922758357796
$ws.Cells.Item(23, 13).Value = 0.5690562345659629
$ws.Cells.Item(23, 15).Value = 2.331714329429843
$ws.Cells.Item(24, 3).Value = 0.06825960667428888
$ws.Cells.Item(24, 4).Value = 0.2251803488729678
$ws.Cells.Item(24, 5).Value = 0.1577196954849107
$ws.Cells.Item(24, 6).Value = 1.073767283960606
$ws.Cells.Item(24, 7).Value = 0.5302274144498824
$ws.Cells.Item(24, 8).Value = 0.6686566082319843
$ws.Cells.Item(24, 10).Value = 0.1538157805362914
$ws.Cells.Item(24, 11).Value = 1.795483521229926
$ws.Cells.Item(24, 13).Value = 0.5089208811288017
$ws.Cells.Item(24, 15).Value = 2.367006240997696
$ws.Cells.Item(25, 3).Value = 0.05716899837759115
$ws.Cells.Item(25, 4).Value = 0.2192114949909012
$ws.Cells.Item(25, 5).Value = 0.1575993397685451
$ws.Cells.Item(25, 6).Value = 1.086324825071799
$ws.Cells.Item(25, 7).Value = 0.5394099294658758
$ws.Cells.Item(25, 8).Value = 0.6813005207687013
$ws.Cells.Item(25, 10).Value = 0.1564641562595419
$ws.Cells.Item(25, 11).Value = 1.490536218305579
$ws.Cells.Item(25, 13).Value = 0.4440676642826205
$ws.Cells.Item(25, 15).Value = 2.412421579825462
